$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Sheet1" to "user data"
$ws.Name = "user data"

# Fix the misspelled "Saftey" -> "Safety" for rows 6-8 in column J
$ws.Range("J6").Value = "Safety"
$ws.Range("J7").Value = "Safety"
$ws.Range("J8").Value = "Safety"

# Update the selected range in the sheet view to J6:J8 with active cell J6
$ws.Range("J6:J8").Select()
